$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41:138 down to 42:139
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new weekly data record
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 44498
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100112017
$ws.Range("G41").Value = "Apio"
$ws.Range("H41").Value = "Americana (o)"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 120
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = 8500
$ws.Range("N41").Value = "`$/docena de matas"
$ws.Range("O41").Value = "Provincia del Elquí"
$ws.Range("P41").Value = 1417
$ws.Range("Q41").Value = 6
$ws.Range("R41").Value = "Hortaliza"
